# Apply updated cryptocurrency price/volume data as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.184.40'
$ws.Range('E2').Value = '  +1.83%  '
$ws.Range('D3').Value = '2.358.71'
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('E5').Value = '  +4.16%  '
$ws.Range('D6').Value = "'243.96"
$ws.Range('E6').Value = '  +3.48%  '
$ws.Range('D7').Value = "'74.46"
$ws.Range('E7').Value = '  +3.17%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').Value = "'0.582"
$ws.Range('E9').Value = '  +25.13%  '
$ws.Range('D10').Value = "'0.103"
$ws.Range('E10').Value = '  +6.12%  '
$ws.Range('E11').Value = '  +18.71%  '
$ws.Range('E12').Value = '  +19.85%  '
$ws.Range('E13').Value = '  +2.03%  '
$ws.Range('D14').Value = '2.710.77'
$ws.Range('E14').Value = '  -0.75%  '
$ws.Range('D15').Value = "'16.91"
$ws.Range('E15').Value = '  +6.38%  '
$ws.Range('D16').Value = "'0.916"
$ws.Range('E16').Value = '  +7.05%  '
$ws.Range('D17').Value = '2.353.51'
$ws.Range('E17').Value = '  -1.09%  '
$ws.Range('D18').Value = '44.366.38'
$ws.Range('E18').Value = '  +2.08%  '
$ws.Range('E19').Value = '  +5.35%  '
$ws.Range('E20').Value = '  +5.81%  '
$ws.Range('D21').Value = "'78.46"
$ws.Range('E21').Value = '  +5.29%  '
$ws.Range('D22').Value = "'256.28"
$ws.Range('E22').Value = '  +2.00%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('E24').Value = '  +4.49%  '
$ws.Range('D25').Value = "'3.76"
$ws.Range('E25').Value = '  -4.81%  '
$ws.Range('D26').Value = "'10.75"
$ws.Range('E26').Value = '  +7.45%  '
$ws.Range('E27').Value = '  +1.26%  '
$ws.Range('D28').Value = "'22.60"
$ws.Range('E28').Value = '  -1.97%  '
$ws.Range('E29').Value = '  +5.26%  '
$ws.Range('D30').Value = "'174.84"
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('E31').Value = '  +3.46%  '
$ws.Range('E32').Value = '  +4.79%  '
$ws.Range('D33').Value = "'5.42"
$ws.Range('E33').Value = '  +8.40%  '
$ws.Range('E34').Value = '  +9.67%  '
$ws.Range('D35').Value = "'5.36"
$ws.Range('E35').Value = '  +5.87%  '
$ws.Range('E36').Value = '  +6.41%  '
$ws.Range('D37').Value = "'2.48"
$ws.Range('E37').Value = '  +0.73%  '
$ws.Range('D38').Value = "'6.58"
$ws.Range('E38').Value = '  -0.53%  '
$ws.Range('E39').Value = '  +7.41%  '
$ws.Range('D40').Value = "'19.48"
$ws.Range('E40').Value = '  +3.89%  '
$ws.Range('D41').Value = "'9.03"
$ws.Range('E41').Value = '  +1.18%  '
$ws.Range('E42').Value = '  +0.36%  '
$ws.Range('D43').Value = "'0.195"
$ws.Range('E43').Value = '  +15.34%  '
$ws.Range('E44').Value = '  +2.82%  '
$ws.Range('E45').Value = '  +11.75%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').Value = "'0.0999"
$ws.Range('E46').Value = '  +5.06%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = "'101.75"
$ws.Range('E47').Value = '  +2.03%  '
$ws.Range('D48').Value = "'1.18"
$ws.Range('E48').Value = '  -0.63%  '
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('D50').Value = '1.457.41'
$ws.Range('E50').Value = '  +0.35%  '
$ws.Range('D51').Value = "'0.000207"
$ws.Range('E51').Value = '  +3.49%  '

# Reset the style on the forced-text cells above so they keep the workbook
# default style (no explicit 's' attribute), same as the original cells.
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
